$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "En tykkää :( 10.4.2017 @ 14:4"
